$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1751412429378531
$ws.Range("C2").Value = 0.5932203389830508
$ws.Range("J2").Value = 0.008474576271186441
$ws.Range("P2").Value = 0.1497175141242938
$ws.Range("S2").Value = 0.07344632768361582
$ws.Range("B3").Value = 0.004464285714285714
$ws.Range("C3").Value = 0.05357142857142857
$ws.Range("J3").Value = 0.02678571428571428
$ws.Range("P3").Value = 0.7366071428571429
$ws.Range("S3").Value = 0.1785714285714286
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("P4").Value = 0.6470588235294118
$ws.Range("S4").Value = 0.3137254901960784
$ws.Range("B6").Value = 0.06140350877192982
$ws.Range("D6").Value = 0.008771929824561403
$ws.Range("F6").Value = 0.04824561403508772
$ws.Range("J6").Value = 0.2675438596491228
$ws.Range("O6").Value = 0.02631578947368421
$ws.Range("Q6").Value = 0.1929824561403509
$ws.Range("R6").Value = 0.04385964912280702
$ws.Range("S6").Value = 0.3508771929824561
$ws.Range("B7").Value = 0.1138211382113821
$ws.Range("D7").Value = 0.02845528455284553
$ws.Range("F7").Value = 0.05691056910569105
$ws.Range("J7").Value = 0.1341463414634146
$ws.Range("O7").Value = 0.01626016260162602
$ws.Range("Q7").Value = 0.1585365853658537
$ws.Range("R7").Value = 0.08130081300813008
$ws.Range("S7").Value = 0.4105691056910569
$ws.Range("B8").Value = 0.1
$ws.Range("D8").Value = 0.0163265306122449
$ws.Range("F8").Value = 0.05306122448979592
$ws.Range("J8").Value = 0.1244897959183673
$ws.Range("O8").Value = 0.0163265306122449
$ws.Range("Q8").Value = 0.2061224489795918
$ws.Range("R8").Value = 0.06122448979591837
$ws.Range("S8").Value = 0.4224489795918367
$ws.Range("B9").Value = 0.1137440758293839
$ws.Range("F9").Value = 0.05687203791469194
$ws.Range("J9").Value = 0.1042654028436019
$ws.Range("Q9").Value = 0.2274881516587678
$ws.Range("R9").Value = 0.0995260663507109
$ws.Range("S9").Value = 0.3981042654028436
$ws.Range("B10").Value = 0.1136212624584718
$ws.Range("D10").Value = 0.02192691029900332
$ws.Range("F10").Value = 0.06378737541528239
$ws.Range("J10").Value = 0.1289036544850498
$ws.Range("O10").Value = 0.01727574750830565
$ws.Range("Q10").Value = 0.2551495016611295
$ws.Range("R10").Value = 0.07109634551495017
$ws.Range("S10").Value = 0.3282392026578073
$ws.Range("G11").Value = 0.1485411140583554
$ws.Range("J11").Value = 0.06631299734748011
$ws.Range("K11").Value = 0.1989389920424403
$ws.Range("L11").Value = 0.5702917771883289
$ws.Range("S11").Value = 0.01591511936339523
$ws.Range("G12").Value = 0.7792792792792793
$ws.Range("J12").Value = 0.1576576576576577
$ws.Range("K12").Value = 0.004504504504504504
$ws.Range("L12").Value = 0.04054054054054054
$ws.Range("S12").Value = 0.01801801801801802
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.5813953488372093
$ws.Range("J13").Value = 0.3488372093023256
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.00796812749003984
$ws.Range("H15").Value = 0.1434262948207171
$ws.Range("I15").Value = 0.0796812749003984
$ws.Range("J15").Value = 0.3745019920318725
$ws.Range("K15").Value = 0.05976095617529881
$ws.Range("M15").Value = 0.01593625498007968
$ws.Range("O15").Value = 0.05577689243027888
$ws.Range("S15").Value = 0.2629482071713147
$ws.Range("F16").Value = 0.01214574898785425
$ws.Range("H16").Value = 0.1700404858299595
$ws.Range("I16").Value = 0.0931174089068826
$ws.Range("J16").Value = 0.4008097165991903
$ws.Range("K16").Value = 0.1255060728744939
$ws.Range("M16").Value = 0.008097165991902834
$ws.Range("N16").Value = 0.004048582995951417
$ws.Range("O16").Value = 0.04453441295546558
$ws.Range("S16").Value = 0.1417004048582996
$ws.Range("F17").Value = 0.008210180623973728
$ws.Range("H17").Value = 0.1412151067323481
$ws.Range("I17").Value = 0.08538587848932677
$ws.Range("J17").Value = 0.4663382594417077
$ws.Range("K17").Value = 0.1001642036124795
$ws.Range("M17").Value = 0.01970443349753695
$ws.Range("O17").Value = 0.0541871921182266
$ws.Range("S17").Value = 0.1247947454844007
$ws.Range("F18").Value = 0.01595744680851064
$ws.Range("H18").Value = 0.1702127659574468
$ws.Range("I18").Value = 0.101063829787234
$ws.Range("J18").Value = 0.4414893617021277
$ws.Range("K18").Value = 0.0797872340425532
$ws.Range("M18").Value = 0.01063829787234043
$ws.Range("O18").Value = 0.05319148936170213
$ws.Range("S18").Value = 0.1276595744680851
$ws.Range("F19").Value = 0.01981599433828733
$ws.Range("H19").Value = 0.2137296532200991
$ws.Range("I19").Value = 0.07077140835102619
$ws.Range("J19").Value = 0.3559801840056617
$ws.Range("K19").Value = 0.1259731068648266
$ws.Range("M19").Value = 0.01840056617126681
$ws.Range("N19").Value = 0.0007077140835102619
$ws.Range("O19").Value = 0.07501769285208776
$ws.Range("S19").Value = 0.1196036801132343
